$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Update status text everywhere: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share this string.)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# --- zh-cn handback: fill in Latest Target File (I), Latest Handback File (J), DateTime (K) ---
$wsZhCn.Range("I2").Value = "656fd5f4-ca5f-4fb5-8028-149023891d81.md"
$wsZhCn.Range("J2").Value = "656fd5f4-ca5f-4fb5-8028-149023891d81.caa9e7a23bf4548d00ce257b07805c1ba501dccb.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-13 16:39:17"

$wsZhCn.Range("I3").Value = "9c488297-dee2-4f1c-94a4-9afcf55844c4.md"
$wsZhCn.Range("J3").Value = "9c488297-dee2-4f1c-94a4-9afcf55844c4.86ec835aa7a4010c29868d94b89ce10663f444ae.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-13 16:39:17"

# --- de-de handback: fill in Latest Target File (I), Latest Handback File (J), DateTime (K) ---
$wsDeDe.Range("I2").Value = "656fd5f4-ca5f-4fb5-8028-149023891d81.md"
$wsDeDe.Range("J2").Value = "656fd5f4-ca5f-4fb5-8028-149023891d81.caa9e7a23bf4548d00ce257b07805c1ba501dccb.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-13 16:39:27"

$wsDeDe.Range("I3").Value = "9c488297-dee2-4f1c-94a4-9afcf55844c4.md"
$wsDeDe.Range("J3").Value = "9c488297-dee2-4f1c-94a4-9afcf55844c4.86ec835aa7a4010c29868d94b89ce10663f444ae.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-13 16:39:27"

# --- Add hyperlinks for the newly-populated "Latest Target File" (I) cells ---
# (Hyperlinks.Add applies the hyperlink cell style automatically.)
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3a24a19b99fb3af36d3803ade0eb15598d547486/e2e/656fd5f4-ca5f-4fb5-8028-149023891d81.md", "", "", "656fd5f4-ca5f-4fb5-8028-149023891d81.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3a24a19b99fb3af36d3803ade0eb15598d547486/e2e/9c488297-dee2-4f1c-94a4-9afcf55844c4.md", "", "", "9c488297-dee2-4f1c-94a4-9afcf55844c4.md") | Out-Null
}

# --- Column widths ---
$wsOverview.Range("E:F").ColumnWidth = 29.9777047293527
$wsZhCn.Range("C:C").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I:J").ColumnWidth = 40
$wsDeDe.Range("C:C").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I:J").ColumnWidth = 40
